# Added some June 2020 info
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 52 (March 2020): add the "Shifted Power" value that was missing
$ws.Range("D52").Value = 164.8

# Row 53 (April 2020)
$ws.Range("D53").Value = 182.17
$ws.Range("E53").Formula = "=D52"
$ws.Range("F53").Value = 65.99
$ws.Range("G53").Formula = "=C53/3"
$ws.Range("H53").Formula = "=(E53+F53)/3"
$ws.Range("J53").Formula = "=G53+H53"
$ws.Range("L53").Formula = "=3 *J53"
$ws.Range("Q53").Value = 458.59
$ws.Range("S53").Formula = "=S52 + L53 - (J53 + M53+ N53 + O53 + P53 + Q53) + I53"
$ws.Range("T53").Formula = "=T52 + J53 - M53"
$ws.Range("X53").Formula = "=X52 + J53 - Q53 + I53"

# Row 54 (May 2020)
$ws.Range("E54").Formula = "=D53"
$ws.Range("F54").Value = 65.99
$ws.Range("G54").Formula = "=C54/3"
$ws.Range("H54").Formula = "=(E54+F54)/3"
$ws.Range("J54").Formula = "=G54+H54"
$ws.Range("L54").Formula = "=3 *J54"
$ws.Range("Q54").Value = 464.39
$ws.Range("S54").Formula = "=S53 + L54 - (J54 + M54+ N54 + O54 + P54 + Q54) + I54"
$ws.Range("T54").Formula = "=T53 + J54 - M54"
$ws.Range("X54").Formula = "=X53 + J54 - Q54 + I54"

# Row 55 (June 2020)
$ws.Range("F55").Value = 65.99
$ws.Range("G55").Formula = "=C55/3"
$ws.Range("H55").Formula = "=(E55+F55)/3"
$ws.Range("J55").Formula = "=G55+H55"
$ws.Range("L55").Formula = "=3 *J55"
$ws.Range("Q55").Value = 403.66
$ws.Range("S55").Formula = "=S54 + L55 - (J55 + M55+ N55 + O55 + P55 + Q55) + I55"
$ws.Range("T55").Formula = "=T54 + J55 - M55"
$ws.Range("X55").Formula = "=X54 + J55 - Q55 + I55"

# Update the frozen-pane selection to reflect the new active cell
$ws.Range("Q56").Select()
